$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 26
$ws.Range("A26").Value = 111273656
$ws.Range("B26").Value = 73696
$ws.Range("E26").Value = 6440
$ws.Range("F26").Value = 'Vitgrynig nållav'
$ws.Range("G26").Value = 'Chaenotheca subroscida'
$ws.Range("H26").Value = '(Eitner) Zahlbr.'
$ws.Range("Q26").Value = 591725.0424782543
$ws.Range("R26").Value = 7043424.7006835

# Row 27
$ws.Range("A27").Value = 111273666
$ws.Range("B27").Value = 96348
$ws.Range("E27").Value = 220787
$ws.Range("F27").Value = 'Knärot'
$ws.Range("G27").Value = 'Goodyera repens'
$ws.Range("H27").Value = '(L.) R. Br.'
$ws.Range("Q27").Value = 591499.5271172373
$ws.Range("R27").Value = 7043317.696102448

# Row 28
$ws.Range("A28").Value = 111273667
$ws.Range("B28").Value = 89423
$ws.Range("E28").Value = 5432
$ws.Range("F28").Value = 'Granticka'
$ws.Range("G28").Value = 'Porodaedalea chrysoloma'
$ws.Range("H28").Value = '(Fr.) Fiasson & Niemelä'
$ws.Range("Q28").Value = 591618.866522243
$ws.Range("R28").Value = 7043352.399297187

# Row 29
$ws.Range("A29").Value = 111273664
$ws.Range("B29").Value = 89590
$ws.Range("D29").Value = 'VU'
$ws.Range("E29").Value = 48
$ws.Range("F29").Value = 'Lappticka'
$ws.Range("G29").Value = 'Amylocystis lapponica'
$ws.Range("H29").Value = '(Romell) Singer'
$ws.Range("Q29").Value = 591673.2841504611
$ws.Range("R29").Value = 7043420.083276978

# Row 30
$ws.Range("A30").Value = 111273655
$ws.Range("B30").Value = 73696
$ws.Range("E30").Value = 6440
$ws.Range("F30").Value = 'Vitgrynig nållav'
$ws.Range("G30").Value = 'Chaenotheca subroscida'
$ws.Range("H30").Value = '(Eitner) Zahlbr.'
$ws.Range("Q30").Value = 591622.4606337334
$ws.Range("R30").Value = 7043398.517451782

# Row 31
$ws.Range("A31").Value = 111273659
$ws.Range("B31").Value = 89845
$ws.Range("D31").Value = 'VU'
$ws.Range("E31").Value = 1209
$ws.Range("F31").Value = 'Rynkskinn'
$ws.Range("G31").Value = 'Phlebia centrifuga'
$ws.Range("H31").Value = 'P.Karst.'
$ws.Range("Q31").Value = 591495.2093399345
$ws.Range("R31").Value = 7043327.847347787

# Row 32
$ws.Range("A32").Value = 111273661
$ws.Range("B32").Value = 89686
$ws.Range("E32").Value = 658
$ws.Range("F32").Value = 'Rosenticka'
$ws.Range("G32").Value = 'Rhodofomes roseus'
$ws.Range("H32").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q32").Value = 591636.9769660851
$ws.Range("R32").Value = 7043422.612332962

# Row 33
$ws.Range("A33").Value = 111273663
$ws.Range("B33").Value = 89686
$ws.Range("D33").Value = 'NT'
$ws.Range("E33").Value = 658
$ws.Range("F33").Value = 'Rosenticka'
$ws.Range("G33").Value = 'Rhodofomes roseus'
$ws.Range("H33").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q33").Value = 591652.4436271309
$ws.Range("R33").Value = 7043413.675855185

# Row 34
$ws.Range("A34").Value = 111273660
$ws.Range("B34").Value = 89845
$ws.Range("D34").Value = 'VU'
$ws.Range("E34").Value = 1209
$ws.Range("F34").Value = 'Rynkskinn'
$ws.Range("G34").Value = 'Phlebia centrifuga'
$ws.Range("H34").Value = 'P.Karst.'
$ws.Range("Q34").Value = 591641.1794901572
$ws.Range("R34").Value = 7043416.478903031

# Row 35
$ws.Range("A35").Value = 111273670
$ws.Range("B35").Value = 77515
$ws.Range("D35").Value = 'NT'
$ws.Range("E35").Value = 6425
$ws.Range("F35").Value = 'Garnlav'
$ws.Range("G35").Value = 'Alectoria sarmentosa'
$ws.Range("H35").Value = '(Ach.) Ach.'
$ws.Range("Q35").Value = 591622.4606337334
$ws.Range("R35").Value = 7043398.517451782

# Row 36
$ws.Range("A36").Value = 111273672
$ws.Range("B36").Value = 77515
$ws.Range("D36").Value = 'NT'
$ws.Range("E36").Value = 6425
$ws.Range("F36").Value = 'Garnlav'
$ws.Range("G36").Value = 'Alectoria sarmentosa'
$ws.Range("H36").Value = '(Ach.) Ach.'
$ws.Range("Q36").Value = 591719.3732997013
$ws.Range("R36").Value = 7043419.6232786

# Row 38
$ws.Range("A38").Value = 111315150
$ws.Range("B38").Value = 89369
$ws.Range("D38").Value = 'LC'
$ws.Range("E38").Value = 5447
$ws.Range("F38").Value = 'Vedticka'
$ws.Range("G38").Value = 'Fuscoporia viticola'
$ws.Range("H38").Value = '(Schwein.) Murrill'
$ws.Range("K38").Value = ""
$ws.Range("Q38").Value = 591671.190636521
$ws.Range("R38").Value = 7043415.108879722

# Row 39
$ws.Range("A39").Value = 111266309
$ws.Range("B39").Value = 77515
$ws.Range("D39").Value = 'NT'
$ws.Range("E39").Value = 6425
$ws.Range("F39").Value = 'Garnlav'
$ws.Range("G39").Value = 'Alectoria sarmentosa'
$ws.Range("H39").Value = '(Ach.) Ach.'
$ws.Range("K39").Value = ""
$ws.Range("Q39").Value = 591747.0822552936
$ws.Range("R39").Value = 7043436.057239689

# Row 40
$ws.Range("A40").Value = 111315151
$ws.Range("B40").Value = 89590
$ws.Range("D40").Value = 'VU'
$ws.Range("E40").Value = 48
$ws.Range("F40").Value = 'Lappticka'
$ws.Range("G40").Value = 'Amylocystis lapponica'
$ws.Range("H40").Value = '(Romell) Singer'
$ws.Range("K40").Value = ""
$ws.Range("Q40").Value = 591670.9593730925
$ws.Range("R40").Value = 7043423.143536596

# Row 42
$ws.Range("A42").Value = 111315146
$ws.Range("K42").Value = ""
$ws.Range("Q42").Value = 591616.7319226691
$ws.Range("R42").Value = 7043364.400079632

# Row 43
$ws.Range("A43").Value = 111315141
$ws.Range("B43").Value = 96348
$ws.Range("D43").Value = 'VU'
$ws.Range("E43").Value = 220787
$ws.Range("F43").Value = 'Knärot'
$ws.Range("G43").Value = 'Goodyera repens'
$ws.Range("H43").Value = '(L.) R. Br.'
$ws.Range("I43").Value = '3'
$ws.Range("Q43").Value = 591486.5005135566
$ws.Range("R43").Value = 7043319.555657836

# Row 44
$ws.Range("A44").Value = 111315149
$ws.Range("Q44").Value = 591670.9593730925
$ws.Range("R44").Value = 7043423.143536596

# Row 45
$ws.Range("A45").Value = 111315142
$ws.Range("B45").Value = 89405
$ws.Range("D45").Value = 'NT'
$ws.Range("E45").Value = 1202
$ws.Range("F45").Value = 'Ullticka'
$ws.Range("G45").Value = 'Phellinidium ferrugineofuscum'
$ws.Range("H45").Value = '(P.Karst.) Fiasson & Niemelä'
$ws.Range("I45").Value = ""
$ws.Range("Q45").Value = 591469.6177441666
$ws.Range("R45").Value = 7043315.49674286

# Row 46
$ws.Range("A46").Value = 111315145
$ws.Range("B46").Value = 89845
$ws.Range("D46").Value = 'VU'
$ws.Range("E46").Value = 1209
$ws.Range("F46").Value = 'Rynkskinn'
$ws.Range("G46").Value = 'Phlebia centrifuga'
$ws.Range("H46").Value = 'P.Karst.'
$ws.Range("Q46").Value = 591478.5830416525
$ws.Range("R46").Value = 7043314.860723522

# Row 47
$ws.Range("A47").Value = 111315148
$ws.Range("Q47").Value = 591645.4590963478
$ws.Range("R47").Value = 7043407.667238996

# Row 48
$ws.Range("A48").Value = 111268460
$ws.Range("B48").Value = 89686
$ws.Range("E48").Value = 658
$ws.Range("F48").Value = 'Rosenticka'
$ws.Range("G48").Value = 'Rhodofomes roseus'
$ws.Range("H48").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("M48").Value = ""

# Row 49
$ws.Range("A49").Value = 111268512
$ws.Range("B49").Value = 56398
$ws.Range("D49").Value = 'NT'
$ws.Range("E49").Value = 100109
$ws.Range("F49").Value = 'Tretåig hackspett'
$ws.Range("G49").Value = 'Picoides tridactylus'
$ws.Range("H49").Value = '(Linnaeus, 1758)'
$ws.Range("I49").Value = ""
$ws.Range("K49").Value = ""
$ws.Range("M49").Value = 'äldre spår'
$ws.Range("AC49").Value = ""
$ws.Range("Q49").Value = 591472.6953434804
$ws.Range("R49").Value = 7043317.372138057

# Row 50
$ws.Range("A50").Value = 111315143
$ws.Range("B50").Value = 89686
$ws.Range("E50").Value = 658
$ws.Range("F50").Value = 'Rosenticka'
$ws.Range("G50").Value = 'Rhodofomes roseus'
$ws.Range("H50").Value = '(Alb. & Schwein.) Kotl. & Pouzar'
$ws.Range("Q50").Value = 591477.5224061215
$ws.Range("R50").Value = 7043320.638036993

# Row 51
$ws.Range("A51").Value = 111315147
$ws.Range("B51").Value = 73696
$ws.Range("E51").Value = 6440
$ws.Range("F51").Value = 'Vitgrynig nållav'
$ws.Range("G51").Value = 'Chaenotheca subroscida'
$ws.Range("H51").Value = '(Eitner) Zahlbr.'
$ws.Range("Q51").Value = 591620.5314988887
$ws.Range("R51").Value = 7043403.376114395

# Row 52
$ws.Range("A52").Value = 111267164
$ws.Range("K52").Value = ""
$ws.Range("Q52").Value = 591635.2558426465
$ws.Range("R52").Value = 7043404.693209249

# Row 53
$ws.Range("A53").Value = 111315139
$ws.Range("B53").Value = 96348
$ws.Range("D53").Value = 'VU'
$ws.Range("E53").Value = 220787
$ws.Range("F53").Value = 'Knärot'
$ws.Range("G53").Value = 'Goodyera repens'
$ws.Range("H53").Value = '(L.) R. Br.'
$ws.Range("I53").Value = '1'
$ws.Range("AC53").Value = 'Plus massor av bladrosetter'
$ws.Range("Q53").Value = 591510.9235177813
$ws.Range("R53").Value = 7043279.155835367
